$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Update existing May values that changed
$ws.Range("B14").Value = 21013.72
$ws.Range("B21").Value = 23465.32

# 2. Insert a new row at position 22 for May day 30 (shifts April/March/Feb down by one row)
$ws.Rows.Item(22).Insert()
$ws.Range("A22").Value = 30
$ws.Range("B22").Value = 3190.35
$ws.Range("C22").Value = 5
$ws.Range("D22").Value = 2025
$ws.Range("E22").Value = "05/2025"

# 3. Delete the February block, which (after the insert above) now lives at rows 62-81
$ws.Range("62:81").Delete()
